$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 11780
$ws.Range("E2").Value = 294
$ws.Range("F2").Value = 294
$ws.Range("G2").Value = -217
$ws.Range("H2").Value = -412
$ws.Range("I2").Value = -392
$ws.Range("J2").Value = -19
$ws.Range("K2").Value = 21579
$ws.Range("L2").Value = 15584
$ws.Range("M2").Value = 5995
$ws.Range("N2").Value = 3712
$ws.Range("O2").Value = 2283
$ws.Range("P2").Value = 1040
$ws.Range("Q2").Value = 1398
$ws.Range("R2").Value = -378
$ws.Range("S2").Value = -839
$ws.Range("T2").Value = 478
$ws.Range("U2").Value = 920
$ws.Range("V2").Value = 13258
$ws.Range("W2").Value = 2.49
$ws.Range("X2").Value = -3.5
$ws.Range("Y2").Value = -10.83
$ws.Range("Z2").Value = -1.86
$ws.Range("AA2").Value = 259.96
$ws.Range("AB2").Value = 254.61
$ws.Range("AC2").Value = -943
$ws.Range("AD2").Value = -2.22
$ws.Range("AE2").Value = 8920
$ws.Range("AF2").Value = 0.23
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 41609310

# Row 3
$ws.Range("D3").Value = 11568
$ws.Range("E3").Value = 724
$ws.Range("F3").Value = 724
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = -44
$ws.Range("J3").Value = 45
$ws.Range("K3").Value = 21703
$ws.Range("L3").Value = 15741
$ws.Range("M3").Value = 5961
$ws.Range("N3").Value = 3831
$ws.Range("O3").Value = 2130
$ws.Range("P3").Value = 1040
$ws.Range("Q3").Value = 807
$ws.Range("R3").Value = -726
$ws.Range("S3").Value = -192
$ws.Range("T3").Value = 771
$ws.Range("U3").Value = 36
$ws.Range("V3").Value = 12702
$ws.Range("W3").Value = 6.26
$ws.Range("X3").Value = 0.01
$ws.Range("Y3").Value = -1.18
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 264.06
$ws.Range("AB3").Value = 255.42
$ws.Range("AC3").Value = -107
$ws.Range("AD3").Value = -28.49
$ws.Range("AE3").Value = 9207
$ws.Range("AF3").Value = 0.33
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 41609310

# Row 4
$ws.Range("D4").Value = 11315
$ws.Range("E4").Value = 421
$ws.Range("F4").Value = 421
$ws.Range("G4").Value = -10
$ws.Range("H4").Value = -14
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = -17
$ws.Range("K4").Value = 19045
$ws.Range("L4").Value = 13130
$ws.Range("M4").Value = 5915
$ws.Range("N4").Value = 3848
$ws.Range("O4").Value = 2067
$ws.Range("P4").Value = 1040
$ws.Range("Q4").Value = 2008
$ws.Range("R4").Value = -187
$ws.Range("S4").Value = -2042
$ws.Range("T4").Value = 190
$ws.Range("U4").Value = 1818
$ws.Range("V4").Value = 10835
$ws.Range("W4").Value = 3.72
$ws.Range("X4").Value = -0.13
$ws.Range("Y4").Value = 0.06
$ws.Range("Z4").Value = -0.07000000000000001
$ws.Range("AA4").Value = 221.98
$ws.Range("AB4").Value = 257.29
$ws.Range("AC4").Value = 5
$ws.Range("AD4").Value = 485.22
$ws.Range("AE4").Value = 9248
$ws.Range("AF4").Value = 0.29
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 41609310

# Row 5
$ws.Range("D5").Value = 10346
$ws.Range("E5").Value = 698
$ws.Range("F5").Value = 698
$ws.Range("G5").Value = 595
$ws.Range("H5").Value = 409
$ws.Range("I5").Value = 267
$ws.Range("J5").Value = 142
$ws.Range("K5").Value = 19157
$ws.Range("L5").Value = 12701
$ws.Range("M5").Value = 6456
$ws.Range("N5").Value = 4105
$ws.Range("O5").Value = 2351
$ws.Range("P5").Value = 1040
$ws.Range("Q5").Value = 1161
$ws.Range("R5").Value = -296
$ws.Range("S5").Value = -525
$ws.Range("T5").Value = 319
$ws.Range("U5").Value = 843
$ws.Range("V5").Value = 10093
$ws.Range("W5").Value = 6.74
$ws.Range("X5").Value = 3.95
$ws.Range("Y5").Value = 6.71
$ws.Range("Z5").Value = 2.14
$ws.Range("AA5").Value = 196.72
$ws.Range("AB5").Value = 283.56
$ws.Range("AC5").Value = 641
$ws.Range("AD5").Value = 3.82
$ws.Range("AE5").Value = 9866
$ws.Range("AF5").Value = 0.25
$ws.Range("AG5").Value = 25
$ws.Range("AH5").Value = 1.02
$ws.Range("AI5").Value = 3.9
$ws.Range("AJ5").Value = 41609310

# Row 6
$ws.Range("D6").Value = 11091
$ws.Range("E6").Value = 1241
$ws.Range("F6").Value = 1241
$ws.Range("G6").Value = 913
$ws.Range("H6").Value = 615
$ws.Range("I6").Value = 347
$ws.Range("K6").Value = 20143
$ws.Range("L6").Value = 13467
$ws.Range("M6").Value = 6676
$ws.Range("N6").Value = 4491
$ws.Range("P6").Value = 1040
$ws.Range("Q6").Value = 855
$ws.Range("R6").Value = -360
$ws.Range("S6").Value = -106
$ws.Range("T6").Value = 354
$ws.Range("U6").Value = 501
$ws.Range("V6").Value = 10323
$ws.Range("W6").Value = 11.19
$ws.Range("X6").Value = 5.55
$ws.Range("Y6").Value = 8.06
$ws.Range("Z6").Value = 3.13
$ws.Range("AA6").Value = 201.74
$ws.Range("AB6").Value = 323.95
$ws.Range("AC6").Value = 833
$ws.Range("AD6").Value = 3.39
$ws.Range("AE6").Value = 10793
$ws.Range("AF6").Value = 0.26
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 1.77
$ws.Range("AI6").Value = 6
$ws.Range("AJ6").Value = 41609310

# Rows 7-9: remove all data beyond column C (keep only A, B, C)
$ws.Range("D7:AJ9").ClearContents()

Write-Output "edit complete"